$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (price / volume) to remain Text-typed cells,
# matching the source data (inline strings), instead of letting Excel
# auto-coerce numeric-looking text into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "89.335.05"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "3.069.63"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "234.48"
$ws.Range("E5").Value = "  +8.13%  "
$ws.Range("D6").Value = "616.61"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "1.05"
$ws.Range("E7").Value = "  -7.21%  "
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "3.069.69"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").Value = "0.704"
$ws.Range("E11").Value = "  -5.99%  "
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "34.79"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "89.210.50"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "5.34"
$ws.Range("E16").Value = "  -6.60%  "
$ws.Range("D17").Value = "3.637.36"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "3.056.09"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").Value = "3.72"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "0.0000211"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "13.70"
$ws.Range("E21").Value = "  -5.85%  "
$ws.Range("D22").Value = "430.02"
$ws.Range("E22").Value = "  -8.54%  "
$ws.Range("D23").Value = "5.37"
$ws.Range("E23").Value = "  +3.74%  "
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("D25").Value = "5.53"
$ws.Range("E25").Value = "  -6.01%  "
$ws.Range("D26").Value = "86.07"
$ws.Range("E26").Value = "  -10.46%  "
$ws.Range("D27").Value = "11.61"
$ws.Range("E27").Value = "  -6.20%  "
$ws.Range("D28").Value = "3.238.97"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +10.09%  "
$ws.Range("D31").Value = "8.99"
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("D32").Value = "0.154"
$ws.Range("E32").Value = "  -6.12%  "
$ws.Range("D33").Value = "0.194"
$ws.Range("E33").Value = "  -10.64%  "
$ws.Range("D34").Value = "25.38"
$ws.Range("E34").Value = "  -6.30%  "
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("D36").Value = "7.03"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").Value = "490.41"
$ws.Range("E37").Value = "  -5.68%  "
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "1.87"
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("B40").Value = "MantraDAO"
$ws.Range("C40").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D40").Value = "3.68"
$ws.Range("E40").Value = "  +52.83%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "1.25"
$ws.Range("E41").Value = "  -6.26%  "
$ws.Range("D42").Value = "0.0884"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").Value = "22.07"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  -9.67%  "
$ws.Range("D46").Value = "155.12"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "1.84"
$ws.Range("E47").Value = "  -7.65%  "
$ws.Range("D48").Value = "0.668"
$ws.Range("E48").Value = "  -8.73%  "
$ws.Range("D49").Value = "44.35"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  -5.68%  "

# Restore the default cell style (the source file has no explicit
# style on these cells) while keeping the Text number format applied
# above so the values stay stored as text.
$ws.Range("D2:E51").Style = "Normal"
